# Update NATMI Ntn3-Neo1 LR-pair TPM results (recomputed with updated TPM
# input data -> ECs "Ligand-expressing cells" count moved from 2/3 to 3/3,
# rippling through detection-rate / expression / specificity / edge-weight
# columns for every row of the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.34281
$ws.Range("H2").Value = 1.02843
$ws.Range("I2").Value = 0.1253250783350286
$ws.Range("J2").Value = 0.1253250783350286
$ws.Range("M2").Value = 3.889188333333333
$ws.Range("N2").Value = 11.667565
$ws.Range("O2").Value = 0.04320782608967438
$ws.Range("P2").Value = 0.04320782608967438
$ws.Range("Q2").Value = 1.33325265255
$ws.Range("R2").Value = 11.99927387295
$ws.Range("S2").Value = 0.005415024189374735
$ws.Range("T2").Value = 0.005415024189374735

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.34281
$ws.Range("H3").Value = 1.02843
$ws.Range("I3").Value = 0.1253250783350286
$ws.Range("J3").Value = 0.1253250783350286
$ws.Range("O3").Value = 0.86451478461177
$ws.Range("P3").Value = 0.8645147846117699
$ws.Range("Q3").Value = 26.67610787361
$ws.Range("R3").Value = 240.08497086249
$ws.Range("S3").Value = 0.1083453831032605
$ws.Range("T3").Value = 0.1083453831032605

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.34281
$ws.Range("H4").Value = 1.02843
$ws.Range("I4").Value = 0.1253250783350286
$ws.Range("J4").Value = 0.1253250783350286
$ws.Range("M4").Value = 0.359731
$ws.Range("N4").Value = 1.079193
$ws.Range("O4").Value = 0.003996513707975397
$ws.Range("P4").Value = 0.003996513707975397
$ws.Range("Q4").Value = 0.12331938411
$ws.Range("R4").Value = 1.10987445699
$ws.Range("S4").Value = 0.0005008633935190324
$ws.Range("T4").Value = 0.0005008633935190324

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.34281
$ws.Range("H5").Value = 1.02843
$ws.Range("I5").Value = 0.1253250783350286
$ws.Range("J5").Value = 0.1253250783350286
$ws.Range("M5").Value = 7.581867
$ws.Range("N5").Value = 22.745601
$ws.Range("O5").Value = 0.08423248315420773
$ws.Range("P5").Value = 0.08423248315420773
$ws.Range("Q5").Value = 2.59913982627
$ws.Range("R5").Value = 23.39225843643
$ws.Range("S5").Value = 0.01055644254965506
$ws.Range("T5").Value = 0.01055644254965506

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.34281
$ws.Range("H6").Value = 1.02843
$ws.Range("I6").Value = 0.1253250783350286
$ws.Range("J6").Value = 0.1253250783350286
$ws.Range("M6").Value = 0.3644006666666667
$ws.Range("N6").Value = 1.093202
$ws.Range("O6").Value = 0.004048392436372474
$ws.Range("P6").Value = 0.004048392436372474
$ws.Range("Q6").Value = 0.12492019254
$ws.Range("R6").Value = 1.12428173286
$ws.Range("S6").Value = 0.0005073650992193177
$ws.Range("T6").Value = 0.0005073650992193177

$ws.Range("I7").Value = 0.5239628476332056
$ws.Range("J7").Value = 0.5239628476332056
$ws.Range("M7").Value = 3.889188333333333
$ws.Range("N7").Value = 11.667565
$ws.Range("O7").Value = 0.04320782608967438
$ws.Range("P7").Value = 0.04320782608967438
$ws.Range("Q7").Value = 5.574102691379444
$ws.Range("R7").Value = 50.166924222415
$ws.Range("S7").Value = 0.0226392955979861
$ws.Range("T7").Value = 0.0226392955979861

$ws.Range("I8").Value = 0.5239628476332056
$ws.Range("J8").Value = 0.5239628476332056
$ws.Range("O8").Value = 0.86451478461177
$ws.Range("P8").Value = 0.8645147846117699
$ws.Range("S8").Value = 0.4529736283661904
$ws.Range("T8").Value = 0.4529736283661903

$ws.Range("I9").Value = 0.5239628476332056
$ws.Range("J9").Value = 0.5239628476332056
$ws.Range("M9").Value = 0.359731
$ws.Range("N9").Value = 1.079193
$ws.Range("O9").Value = 0.003996513707975397
$ws.Range("P9").Value = 0.003996513707975397
$ws.Range("Q9").Value = 0.5155773810403333
$ws.Range("R9").Value = 4.640196429363001
$ws.Range("S9").Value = 0.002094024703035931
$ws.Range("T9").Value = 0.002094024703035931

$ws.Range("I10").Value = 0.5239628476332056
$ws.Range("J10").Value = 0.5239628476332056
$ws.Range("M10").Value = 7.581867
$ws.Range("N10").Value = 22.745601
$ws.Range("O10").Value = 0.08423248315420773
$ws.Range("P10").Value = 0.08423248315420773
$ws.Range("Q10").Value = 10.866561767699
$ws.Range("R10").Value = 97.79905590929101
$ws.Range("S10").Value = 0.0441346917366947
$ws.Range("T10").Value = 0.0441346917366947

$ws.Range("I11").Value = 0.5239628476332056
$ws.Range("J11").Value = 0.5239628476332056
$ws.Range("M11").Value = 0.3644006666666667
$ws.Range("N11").Value = 1.093202
$ws.Range("O11").Value = 0.004048392436372474
$ws.Range("P11").Value = 0.004048392436372474
$ws.Range("Q11").Value = 0.5222700889535555
$ws.Range("R11").Value = 4.700430800582001
$ws.Range("S11").Value = 0.002121207229298453
$ws.Range("T11").Value = 0.002121207229298453

$ws.Range("G12").Value = 0.485713
$ws.Range("H12").Value = 1.457139
$ws.Range("I12").Value = 0.1775678065789847
$ws.Range("J12").Value = 0.1775678065789847
$ws.Range("M12").Value = 3.889188333333333
$ws.Range("N12").Value = 11.667565
$ws.Range("O12").Value = 0.04320782608967438
$ws.Range("P12").Value = 0.04320782608967438
$ws.Range("Q12").Value = 1.889029332948333
$ws.Range("R12").Value = 17.001263996535
$ws.Range("S12").Value = 0.007672318905789709
$ws.Range("T12").Value = 0.007672318905789709

$ws.Range("G13").Value = 0.485713
$ws.Range("H13").Value = 1.457139
$ws.Range("I13").Value = 0.1775678065789847
$ws.Range("J13").Value = 0.1775678065789847
$ws.Range("O13").Value = 0.86451478461177
$ws.Range("P13").Value = 0.8645147846117699
$ws.Range("Q13").Value = 37.79624976988633
$ws.Range("R13").Value = 340.166247928977
$ws.Range("S13").Value = 0.1535099940586154
$ws.Range("T13").Value = 0.1535099940586154

$ws.Range("G14").Value = 0.485713
$ws.Range("H14").Value = 1.457139
$ws.Range("I14").Value = 0.1775678065789847
$ws.Range("J14").Value = 0.1775678065789847
$ws.Range("M14").Value = 0.359731
$ws.Range("N14").Value = 1.079193
$ws.Range("O14").Value = 0.003996513707975397
$ws.Range("P14").Value = 0.003996513707975397
$ws.Range("Q14").Value = 0.174726023203
$ws.Range("R14").Value = 1.572534208827
$ws.Range("S14").Value = 0.0007096521730880364
$ws.Range("T14").Value = 0.0007096521730880364

$ws.Range("G15").Value = 0.485713
$ws.Range("H15").Value = 1.457139
$ws.Range("I15").Value = 0.1775678065789847
$ws.Range("J15").Value = 0.1775678065789847
$ws.Range("M15").Value = 7.581867
$ws.Range("N15").Value = 22.745601
$ws.Range("O15").Value = 0.08423248315420773
$ws.Range("P15").Value = 0.08423248315420773
$ws.Range("Q15").Value = 3.682611366171
$ws.Range("R15").Value = 33.143502295539
$ws.Range("S15").Value = 0.01495697727639395
$ws.Range("T15").Value = 0.01495697727639395

$ws.Range("G16").Value = 0.485713
$ws.Range("H16").Value = 1.457139
$ws.Range("I16").Value = 0.1775678065789847
$ws.Range("J16").Value = 0.1775678065789847
$ws.Range("M16").Value = 0.3644006666666667
$ws.Range("N16").Value = 1.093202
$ws.Range("O16").Value = 0.004048392436372474
$ws.Range("P16").Value = 0.004048392436372474
$ws.Range("Q16").Value = 0.1769941410086667
$ws.Range("R16").Value = 1.592947269078
$ws.Range("S16").Value = 0.0007188641650976122
$ws.Range("T16").Value = 0.0007188641650976122

$ws.Range("G17").Value = 0.2505783333333333
$ws.Range("H17").Value = 0.751735
$ws.Range("I17").Value = 0.09160686460156038
$ws.Range("J17").Value = 0.09160686460156038
$ws.Range("M17").Value = 3.889188333333333
$ws.Range("N17").Value = 11.667565
$ws.Range("O17").Value = 0.04320782608967438
$ws.Range("P17").Value = 0.04320782608967438
$ws.Range("Q17").Value = 0.9745463305861112
$ws.Range("R17").Value = 8.770916975275
$ws.Range("S17").Value = 0.003958133474324569
$ws.Range("T17").Value = 0.003958133474324569

$ws.Range("G18").Value = 0.2505783333333333
$ws.Range("H18").Value = 0.751735
$ws.Range("I18").Value = 0.09160686460156038
$ws.Range("J18").Value = 0.09160686460156038
$ws.Range("O18").Value = 0.86451478461177
$ws.Range("P18").Value = 0.8645147846117699
$ws.Range("Q18").Value = 19.49900717828944
$ws.Range("R18").Value = 175.491064604605
$ws.Range("S18").Value = 0.07919548881997755
$ws.Range("T18").Value = 0.07919548881997754

$ws.Range("G19").Value = 0.2505783333333333
$ws.Range("H19").Value = 0.751735
$ws.Range("I19").Value = 0.09160686460156038
$ws.Range("J19").Value = 0.09160686460156038
$ws.Range("M19").Value = 0.359731
$ws.Range("N19").Value = 1.079193
$ws.Range("O19").Value = 0.003996513707975397
$ws.Range("P19").Value = 0.003996513707975397
$ws.Range("Q19").Value = 0.09014079442833334
$ws.Range("R19").Value = 0.8112671498550001
$ws.Range("S19").Value = 0.0003661080901247822
$ws.Range("T19").Value = 0.0003661080901247822

$ws.Range("G20").Value = 0.2505783333333333
$ws.Range("H20").Value = 0.751735
$ws.Range("I20").Value = 0.09160686460156038
$ws.Range("J20").Value = 0.09160686460156038
$ws.Range("M20").Value = 7.581867
$ws.Range("N20").Value = 22.745601
$ws.Range("O20").Value = 0.08423248315420773
$ws.Range("P20").Value = 0.08423248315420773
$ws.Range("Q20").Value = 1.899851596415
$ws.Range("R20").Value = 17.098664367735
$ws.Range("S20").Value = 0.007716273679360723
$ws.Range("T20").Value = 0.007716273679360723

$ws.Range("G21").Value = 0.2505783333333333
$ws.Range("H21").Value = 0.751735
$ws.Range("I21").Value = 0.09160686460156038
$ws.Range("J21").Value = 0.09160686460156038
$ws.Range("M21").Value = 0.3644006666666667
$ws.Range("N21").Value = 1.093202
$ws.Range("O21").Value = 0.004048392436372474
$ws.Range("P21").Value = 0.004048392436372474
$ws.Range("Q21").Value = 0.09131091171888889
$ws.Range("R21").Value = 0.82179820547
$ws.Range("S21").Value = 0.0003708605377727544
$ws.Range("T21").Value = 0.0003708605377727544

$ws.Range("G22").Value = 0.2230346666666667
$ws.Range("H22").Value = 0.669104
$ws.Range("I22").Value = 0.08153740285122078
$ws.Range("J22").Value = 0.08153740285122078
$ws.Range("M22").Value = 3.889188333333333
$ws.Range("N22").Value = 11.667565
$ws.Range("O22").Value = 0.04320782608967438
$ws.Range("P22").Value = 0.04320782608967438
$ws.Range("Q22").Value = 0.8674238235288889
$ws.Range("R22").Value = 7.80681441176
$ws.Range("S22").Value = 0.003523053922199267
$ws.Range("T22").Value = 0.003523053922199267

$ws.Range("G23").Value = 0.2230346666666667
$ws.Range("H23").Value = 0.669104
$ws.Range("I23").Value = 0.08153740285122078
$ws.Range("J23").Value = 0.08153740285122078
$ws.Range("O23").Value = 0.86451478461177
$ws.Range("P23").Value = 0.8645147846117699
$ws.Range("Q23").Value = 17.35566881816355
$ws.Range("R23").Value = 156.201019363472
$ws.Range("S23").Value = 0.07049029026372626
$ws.Range("T23").Value = 0.07049029026372625

$ws.Range("G24").Value = 0.2230346666666667
$ws.Range("H24").Value = 0.669104
$ws.Range("I24").Value = 0.08153740285122078
$ws.Range("J24").Value = 0.08153740285122078
$ws.Range("M24").Value = 0.359731
$ws.Range("N24").Value = 1.079193
$ws.Range("O24").Value = 0.003996513707975397
$ws.Range("P24").Value = 0.003996513707975397
$ws.Range("Q24").Value = 0.08023248367466668
$ws.Range("R24").Value = 0.722092353072
$ws.Range("S24").Value = 0.0003258653482076161
$ws.Range("T24").Value = 0.0003258653482076161

$ws.Range("G25").Value = 0.2230346666666667
$ws.Range("H25").Value = 0.669104
$ws.Range("I25").Value = 0.08153740285122078
$ws.Range("J25").Value = 0.08153740285122078
$ws.Range("M25").Value = 7.581867
$ws.Range("N25").Value = 22.745601
$ws.Range("O25").Value = 0.08423248315420773
$ws.Range("P25").Value = 0.08423248315420773
$ws.Range("Q25").Value = 1.691019179056
$ws.Range("R25").Value = 15.219172611504
$ws.Range("S25").Value = 0.006868097912103304
$ws.Range("T25").Value = 0.006868097912103304

$ws.Range("G26").Value = 0.2230346666666667
$ws.Range("H26").Value = 0.669104
$ws.Range("I26").Value = 0.08153740285122078
$ws.Range("J26").Value = 0.08153740285122078
$ws.Range("M26").Value = 0.3644006666666667
$ws.Range("N26").Value = 1.093202
$ws.Range("O26").Value = 0.004048392436372474
$ws.Range("P26").Value = 0.004048392436372474
$ws.Range("Q26").Value = 0.08127398122311111
$ws.Range("R26").Value = 0.731465831008
$ws.Range("S26").Value = 0.0003300954049843377
$ws.Range("T26").Value = 0.0003300954049843377
